$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four rows (12:15) that held the tasks for the old week (41);
# everything below shifts up.
$ws.Rows("12:15").Delete()

# Week number: 41 -> 46 for every data row (2:11)
$ws.Range("B2:B11").Value = 46

# New task text for the 5 "days" being added this week (rows 7:11)
$ws.Range("C7").Value = "Ubah promo Happy Hour di MOKA Pos Complete Me"
$ws.Range("C8").Value = "Revisi Harga online di Onlien Store Complete Me & Hyangyu"
$ws.Range("C9").Value = "Maintenance Komputer Server & CCTV Toko Mas an an"
$ws.Range("C10").Value = "Maintenance CCTV,Internet & Printer admin Complete Me"
$ws.Range("C11").Value = "Maintenace internet & CCTV Amazy & Central Kitchen"

# Column C needs to be widened to fit the longer task descriptions
$ws.Columns("C").ColumnWidth = 113.28515625

# Selection moves back to the top of the task column
$ws.Range("C2").Select()
